$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '97.570.88'
$c.Style = $s
$ws.Range('E2').Value = '  -1.55%  '

$c = $ws.Range('D3')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.401.10'
$c.Style = $s
$ws.Range('E3').Value = '  +3.44%  '

$c = $ws.Range('D4')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = $s
$ws.Range('E4').Value = '  -0.04%  '

$c = $ws.Range('D5')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '253.91'
$c.Style = $s
$ws.Range('E5').Value = '  -0.33%  '

$c = $ws.Range('D6')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '648.02'
$c.Style = $s
$ws.Range('E6').Value = '  +3.84%  '

$ws.Range('E7').Value = '  -0.80%  '

$c = $ws.Range('D8')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.423'
$c.Style = $s
$ws.Range('E8').Value = '  +5.40%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range('D9')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.04'
$c.Style = $s
$ws.Range('E9').Value = '  +5.67%  '

$ws.Range('B10').Value = 'USDC'
$ws.Range('C10').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$c = $ws.Range('D10')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = $s
$ws.Range('E10').Value = '  -0.02%  '

$c = $ws.Range('D11')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.398.38'
$c.Style = $s
$ws.Range('E11').Value = '  +3.44%  '

$ws.Range('E12').Value = '  +4.67%  '

$c = $ws.Range('D13')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '41.09'
$c.Style = $s
$ws.Range('E13').Value = '  +3.24%  '

$c = $ws.Range('D14')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.23'
$c.Style = $s
$ws.Range('E14').Value = '  +13.36%  '

$ws.Range('E15').Value = '  +2.64%  '

$c = $ws.Range('D16')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '97.288.95'
$c.Style = $s
$ws.Range('E16').Value = '  -1.52%  '

$c = $ws.Range('D17')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.028.96'
$c.Style = $s
$ws.Range('E17').Value = '  +3.35%  '

$c = $ws.Range('D18')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.40'
$c.Style = $s
$ws.Range('E18').Value = '  +32.02%  '

$c = $ws.Range('D19')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.397.06'
$c.Style = $s
$ws.Range('E19').Value = '  +3.15%  '

$c = $ws.Range('D20')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '17.30'
$c.Style = $s
$ws.Range('E20').Value = '  +13.28%  '

$ws.Range('E21').Value = '  +14.48%  '

$ws.Range('B22').Value = 'Stellar'
$ws.Range('C22').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D22')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.479'
$c.Style = $s
$ws.Range('E22').Value = '  +38.57%  '

$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range('D23')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.38'
$c.Style = $s
$ws.Range('E23').Value = '  -2.46%  '

$c = $ws.Range('D24')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '500.69'
$c.Style = $s
$ws.Range('E24').Value = '  +2.15%  '

$c = $ws.Range('D25')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0000203'
$c.Style = $s
$ws.Range('E25').Value = '  +0.17%  '

$c = $ws.Range('D26')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.03'
$c.Style = $s
$ws.Range('E26').Value = '  +6.64%  '

$c = $ws.Range('D27')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '95.38'
$c.Style = $s
$ws.Range('E27').Value = '  +7.01%  '

$c = $ws.Range('D28')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '12.55'
$c.Style = $s
$ws.Range('E28').Value = '  +3.56%  '

$c = $ws.Range('D29')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.574.80'
$c.Style = $s
$ws.Range('E29').Value = '  +3.26%  '

$ws.Range('E30').Value = '  +9.83%  '

$c = $ws.Range('D31')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.197'
$c.Style = $s
$ws.Range('E31').Value = '  +4.44%  '

$ws.Range('E32').Value = '  -0.30%  '

$c = $ws.Range('D33')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '11.19'
$c.Style = $s
$ws.Range('E33').Value = '  +6.81%  '

$c = $ws.Range('D34')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = $s
$ws.Range('E34').Value = '  -0.16%  '

$ws.Range('E35').Value = '  +17.93%  '

$c = $ws.Range('D36')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '29.42'
$c.Style = $s
$ws.Range('E36').Value = '  +4.87%  '

$c = $ws.Range('D37')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.23'
$c.Style = $s
$ws.Range('E37').Value = '  +14.34%  '

$c = $ws.Range('D38')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.63'
$c.Style = $s
$ws.Range('E38').Value = '  +5.16%  '

$ws.Range('E39').Value = '  +1.78%  '

$ws.Range('E40').Value = '  +12.38%  '

$c = $ws.Range('D41')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '505.65'
$c.Style = $s
$ws.Range('E41').Value = '  +3.57%  '

$c = $ws.Range('D42')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '24.70'
$c.Style = $s
$ws.Range('E42').Value = '  -0.23%  '

$c = $ws.Range('D43')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.853'
$c.Style = $s
$ws.Range('E43').Value = '  +10.06%  '

$c = $ws.Range('D44')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '3.65'
$c.Style = $s
$ws.Range('E44').Value = '  -2.79%  '

$c = $ws.Range('D45')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0411'
$c.Style = $s
$ws.Range('E45').Value = '  +20.69%  '

$c = $ws.Range('D46')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.47'
$c.Style = $s
$ws.Range('E46').Value = '  +14.46%  '

$ws.Range('E47').Value = '  +0.01%  '

$ws.Range('E48').Value = '  +2.46%  '

$c = $ws.Range('D49')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.12'
$c.Style = $s
$ws.Range('E49').Value = '  +10.65%  '

$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D50')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.55'
$c.Style = $s
$ws.Range('E50').Value = '  +13.32%  '

$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D51')
$s = $c.Style
$c.NumberFormat = '@'
$c.Value = '51.44'
$c.Style = $s
$ws.Range('E51').Value = '  +11.17%  '
